$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.084.28'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '2.371.60'
$ws.Range("E3").Value = '  +1.77%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.03'
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.40'
$ws.Range("E6").Value = '  +1.04%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.480'
$ws.Range("E9").Value = '  -2.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.29'
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.124'
$ws.Range("E11").Value = '  +2.46%  '
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.37'
$ws.Range("E13").Value = '  -1.95%  '
$ws.Range("E14").Value = '  +0.64%  '
$ws.Range("D15").Value = '2.739.12'
$ws.Range("E15").Value = '  +1.95%  '
$ws.Range("D16").Value = '2.370.99'
$ws.Range("E16").Value = '  +2.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.797'
$ws.Range("E17").Value = '  +0.75%  '
$ws.Range("D18").Value = '43.123.05'
$ws.Range("E18").Value = '  +0.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.96'
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.27'
$ws.Range("E20").Value = '  +1.33%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.95'
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.64'
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.20'
$ws.Range("E24").Value = '  -1.05%  '
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  +7.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.34'
$ws.Range("E29").Value = '  +2.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.23'
$ws.Range("E30").Value = '  +3.00%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  +0.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.60'
$ws.Range("E33").Value = '  +2.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.110'
$ws.Range("E34").Value = '  +9.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0731'
$ws.Range("E35").Value = '  -3.13%  '
$ws.Range("E36").Value = '  +1.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '127.04'
$ws.Range("E37").Value = '  +2.41%  '
$ws.Range("E38").Value = '  +4.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.31'
$ws.Range("E39").Value = '  -1.56%  '
$ws.Range("E40").Value = '  -2.31%  '
$ws.Range("E41").Value = '  -0.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.72'
$ws.Range("E42").Value = '  -5.93%  '
$ws.Range("D43").Value = '1.929.05'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0278'
$ws.Range("E44").Value = '  -1.29%  '
$ws.Range("E45").Value = '  +2.66%  '
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.20'
$ws.Range("E47").Value = '  -8.59%  '
$ws.Range("D48").Value = '2.598.40'
$ws.Range("E48").Value = '  +1.79%  '
$ws.Range("E49").Value = '  +2.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '71.41'
$ws.Range("E50").Value = '  -0.51%  '
$ws.Range("E51").Value = '  +1.25%  '
